# Add a new product row into the Maestro article list.
# A brand-new row is inserted at row 3 (pushing all the existing rows
# below it down by one, rows 2..58 -> 3..59), and the freed-up row 3 is
# populated with the data for the new "Armonia" milk product.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (2..58) down by one, to rows 3..59,
# carrying their formatting along with them.
$ws.Range("A2:P58").Copy()
$ws.Range("A3:P59").PasteSpecial()
$excel.CutCopyMode = $false

# The brand-new last row (59, beyond the sheet's previous used range)
# doesn't pick up column A's number format from the paste above, so
# reapply it explicitly to keep the barcode formatted like every other row.
$ws.Cells.Item(59, 1).NumberFormat = $ws.Cells.Item(58, 1).NumberFormat

# Populate the newly freed-up row 3 with the new article's data.
$ws.Cells.Item(3, 1).Value = 7790742336200      # A: Codigo
$ws.Cells.Item(3, 2).Value = "Leche"            # B: TipoArtículo
$ws.Cells.Item(3, 3).Value = "larga vida"       # C: Marca
$ws.Cells.Item(3, 4).Value = "parcialmente descremada 2%"  # D: Descripción
$ws.Cells.Item(3, 5).Value = "Armonia"          # E: Variedad
$ws.Cells.Item(3, 6).Value = 1                  # F: ContenidoNeto
$ws.Cells.Item(3, 7).Value = "lt."              # G: UnidadDeMedida
$ws.Cells.Item(3, 8).Value = "tetra brik"       # H: Packaging
$ws.Cells.Item(3, 9).Value = "Leches"           # I: Familia
$ws.Cells.Item(3, 10).Value = "Argentina"       # J: PaisOrigen
$ws.Cells.Item(3, 11).Value = 12                # K: UnidadesPorBulto
$ws.Cells.Item(3, 12).Value = $false            # L: Pesable
$ws.Cells.Item(3, 13).Value = $true             # M: TieneVencimiento
$ws.Cells.Item(3, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790742336200.png"  # N: Imagen
$ws.Cells.Item(3, 15).Value = $true             # O: ImagenExactaDelArticulo
$ws.Cells.Item(3, 16).Value = $true             # P: DescripciónConPackaging

# Column D now holds a longer description than before, so its best-fit
# width grows to accommodate it.
$ws.Columns.Item(4).AutoFit()
